$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the value_1 / value_2 numbers for rows 2 and 3, and point both
# row's "result" formula at row 3 (A3/B3) instead of row 2 (A2/B2).
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 0
$ws.Range("C2").Formula = "=A3/B3"

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("C3").Formula = "=A3/B3"

# Move the active selection to E8.
$ws.Range("E8").Select()
